$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data range (including new rows) is formatted as text so AIH/CBO/CNS
# codes are preserved as strings (matching the original inlineStr content),
# not reinterpreted as numbers.
$ws.Range("A2:C30").NumberFormat = "@"

# Update existing rows 2-19 and add new rows 20-30 with the new AIH results data
$ws.Cells.Item(2, 1).Value = "312110191680-1"
$ws.Cells.Item(2, 2).Value = "225270"
$ws.Cells.Item(2, 3).Value = "705205440440872"
$ws.Cells.Item(3, 1).Value = "312110192110-2"
$ws.Cells.Item(3, 2).Value = "225270"
$ws.Cells.Item(3, 3).Value = "705205440440872"
$ws.Cells.Item(4, 1).Value = "312110192122-3"
$ws.Cells.Item(4, 2).Value = "225270"
$ws.Cells.Item(4, 3).Value = "705205440440872"
$ws.Cells.Item(5, 1).Value = "312110192122-3"
$ws.Cells.Item(5, 2).Value = "225270"
$ws.Cells.Item(5, 3).Value = "705205440440872"
$ws.Cells.Item(6, 1).Value = "312110192147-6"
$ws.Cells.Item(6, 2).Value = "225270"
$ws.Cells.Item(6, 3).Value = "705205440440872"
$ws.Cells.Item(7, 1).Value = "312110192175-1"
$ws.Cells.Item(7, 2).Value = "225270"
$ws.Cells.Item(7, 3).Value = "127980992280001"
$ws.Cells.Item(8, 1).Value = "312110192175-1"
$ws.Cells.Item(8, 2).Value = "225270"
$ws.Cells.Item(8, 3).Value = "127980992280001"
$ws.Cells.Item(9, 1).Value = "312110192175-1"
$ws.Cells.Item(9, 2).Value = "225270"
$ws.Cells.Item(9, 3).Value = "127980992280001"
$ws.Cells.Item(10, 1).Value = "312110192185-0"
$ws.Cells.Item(10, 2).Value = "225270"
$ws.Cells.Item(10, 3).Value = "127980992280001"
$ws.Cells.Item(11, 1).Value = "312110192185-0"
$ws.Cells.Item(11, 2).Value = "225270"
$ws.Cells.Item(11, 3).Value = "127980992280001"
$ws.Cells.Item(12, 1).Value = "312110192185-0"
$ws.Cells.Item(12, 2).Value = "225270"
$ws.Cells.Item(12, 3).Value = "127980992280001"
$ws.Cells.Item(13, 1).Value = "312110192204-8"
$ws.Cells.Item(13, 2).Value = "225270"
$ws.Cells.Item(13, 3).Value = "190149628110005"
$ws.Cells.Item(14, 1).Value = "312110192204-8"
$ws.Cells.Item(14, 2).Value = "225270"
$ws.Cells.Item(14, 3).Value = "190149628110005"
$ws.Cells.Item(15, 1).Value = "312110192900-0"
$ws.Cells.Item(15, 2).Value = "225270"
$ws.Cells.Item(15, 3).Value = "705205440440872"
$ws.Cells.Item(16, 1).Value = "312110192901-1"
$ws.Cells.Item(16, 2).Value = "225270"
$ws.Cells.Item(16, 3).Value = "705205440440872"
$ws.Cells.Item(17, 1).Value = "312110193161-8"
$ws.Cells.Item(17, 2).Value = "225270"
$ws.Cells.Item(17, 3).Value = "204322025140005"
$ws.Cells.Item(18, 1).Value = "312110193161-8"
$ws.Cells.Item(18, 2).Value = "225270"
$ws.Cells.Item(18, 3).Value = "204322025140005"
$ws.Cells.Item(19, 1).Value = "312110193483-0"
$ws.Cells.Item(19, 2).Value = "225270"
$ws.Cells.Item(19, 3).Value = "980016286834678"
$ws.Cells.Item(20, 1).Value = "312110193483-0"
$ws.Cells.Item(20, 2).Value = "225270"
$ws.Cells.Item(20, 3).Value = "980016286834678"
$ws.Cells.Item(21, 1).Value = "312110193483-0"
$ws.Cells.Item(21, 2).Value = "225270"
$ws.Cells.Item(21, 3).Value = "980016286834678"
$ws.Cells.Item(22, 1).Value = "312110194339-9"
$ws.Cells.Item(22, 2).Value = "225270"
$ws.Cells.Item(22, 3).Value = "204322025140005"
$ws.Cells.Item(23, 1).Value = "312110194339-9"
$ws.Cells.Item(23, 2).Value = "225270"
$ws.Cells.Item(23, 3).Value = "204322025140005"
$ws.Cells.Item(24, 1).Value = "312110194953-7"
$ws.Cells.Item(24, 2).Value = "225270"
$ws.Cells.Item(24, 3).Value = "204322025140005"
$ws.Cells.Item(25, 1).Value = "312110194953-7"
$ws.Cells.Item(25, 2).Value = "225270"
$ws.Cells.Item(25, 3).Value = "204322025140005"
$ws.Cells.Item(26, 1).Value = "312150252863-9"
$ws.Cells.Item(26, 2).Value = "225270"
$ws.Cells.Item(26, 3).Value = "980016286834678"
$ws.Cells.Item(27, 1).Value = "312150252864-0"
$ws.Cells.Item(27, 2).Value = "225270"
$ws.Cells.Item(27, 3).Value = "980016286834678"
$ws.Cells.Item(28, 1).Value = "312150252865-0"
$ws.Cells.Item(28, 2).Value = "225270"
$ws.Cells.Item(28, 3).Value = "980016286834678"
$ws.Cells.Item(29, 1).Value = "312150252871-6"
$ws.Cells.Item(29, 2).Value = "225270"
$ws.Cells.Item(29, 3).Value = "980016286834678"
$ws.Cells.Item(30, 1).Value = "312150252874-9"
$ws.Cells.Item(30, 2).Value = "225270"
$ws.Cells.Item(30, 3).Value = "980016286834678"

Write-Host "Updated rows 2-30 with new AIH results data"
